# Generate Report for Handback
# The handback report now reflects that 02df0fce-b5c3-456a-ae67-c5db7d8aab18.md
# has been handed back (in sync with en-US), while b7730ea5-0d78-48b3-8419-1bed7197d7b1.md
# keeps its previously-handed-back data. Row 2 / Row 3 swap their file identity on
# every sheet and the zh-cn / de-de detail sheets get refreshed handoff/handback
# timestamps plus the resolved (no more stale-version) status for 02df0fce.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "02df0fce-b5c3-456a-ae67-c5db7d8aab18.md"
$ws.Range("B2").Value = "e2e\02df0fce-b5c3-456a-ae67-c5db7d8aab18.md"
$ws.Range("G2").Value = "2016-08-13 10:59:35"

$ws.Range("A3").Value = "b7730ea5-0d78-48b3-8419-1bed7197d7b1.md"
$ws.Range("B3").Value = "e2e\b7730ea5-0d78-48b3-8419-1bed7197d7b1.md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-08-13 10:58:29"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "02df0fce-b5c3-456a-ae67-c5db7d8aab18.md"
$ws.Range("G2").Value = "02df0fce-b5c3-456a-ae67-c5db7d8aab18.7bf2e401af40a4b86c91b34c0e9ef4013b37224d.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-13 10:59:27"
$ws.Range("I2").Value = "02df0fce-b5c3-456a-ae67-c5db7d8aab18.md"
$ws.Range("J2").Value = "02df0fce-b5c3-456a-ae67-c5db7d8aab18.7bf2e401af40a4b86c91b34c0e9ef4013b37224d.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-13 10:59:55"

$ws.Range("A3").Value = "b7730ea5-0d78-48b3-8419-1bed7197d7b1.md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "b7730ea5-0d78-48b3-8419-1bed7197d7b1.95755712bd1477a3ddf9adbb182fe073aceeb6bb.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-13 10:58:22"
$ws.Range("I3").Value = "b7730ea5-0d78-48b3-8419-1bed7197d7b1.md"
$ws.Range("J3").Value = "b7730ea5-0d78-48b3-8419-1bed7197d7b1.95755712bd1477a3ddf9adbb182fe073aceeb6bb.zh-cn.xlf"
$ws.Range("P3").Value = ""

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "02df0fce-b5c3-456a-ae67-c5db7d8aab18.md"
$ws.Range("G2").Value = "02df0fce-b5c3-456a-ae67-c5db7d8aab18.7bf2e401af40a4b86c91b34c0e9ef4013b37224d.de-de.xlf"
$ws.Range("H2").Value = "2016-08-13 10:59:35"
$ws.Range("I2").Value = "02df0fce-b5c3-456a-ae67-c5db7d8aab18.md"
$ws.Range("J2").Value = "02df0fce-b5c3-456a-ae67-c5db7d8aab18.7bf2e401af40a4b86c91b34c0e9ef4013b37224d.de-de.xlf"
$ws.Range("K2").Value = "2016-08-13 11:00:12"

$ws.Range("A3").Value = "b7730ea5-0d78-48b3-8419-1bed7197d7b1.md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "b7730ea5-0d78-48b3-8419-1bed7197d7b1.95755712bd1477a3ddf9adbb182fe073aceeb6bb.de-de.xlf"
$ws.Range("H3").Value = "2016-08-13 10:58:29"
$ws.Range("I3").Value = "b7730ea5-0d78-48b3-8419-1bed7197d7b1.md"
$ws.Range("J3").Value = "b7730ea5-0d78-48b3-8419-1bed7197d7b1.95755712bd1477a3ddf9adbb182fe073aceeb6bb.de-de.xlf"
$ws.Range("P3").Value = ""
